# Update KNN/data.xlsx worksheet:
#  - rename header "Дом" -> "Дом_далеко" (E1)
#  - rename value "Не доволен" -> "Недоволен" wherever it occurs (H3,H4,H8,H9,H11)
#  - apply a "0.0" number format to D11 (existing half-point value) and to the
#    newly added D12 half-point value
#  - append two new observation rows (12 and 13)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- header rename ---------------------------------------------------------
$ws.Range("E1").Value = "Дом_далеко"

# --- fix wording of the "dissatisfied" label on existing rows --------------
$ws.Range("H3").Value = "Недоволен"
$ws.Range("H4").Value = "Недоволен"
$ws.Range("H8").Value = "Недоволен"
$ws.Range("H9").Value = "Недоволен"
$ws.Range("H11").Value = "Недоволен"

# --- give the half-point "Подъем" values a one-decimal number format -------
$ws.Range("D11").NumberFormat = "0.0"

# --- append new row 12 -------------------------------------------------------
$ws.Cells.Item(12, 1).Value = 22
$ws.Cells.Item(12, 2).Value = "м"
$ws.Cells.Item(12, 3).Value = "Сова"
$ws.Cells.Item(12, 4).Value = 8.5
$ws.Cells.Item(12, 4).NumberFormat = "0.0"
$ws.Cells.Item(12, 5).Value = "Да"
$ws.Cells.Item(12, 6).Value = "Экстраверт"
$ws.Cells.Item(12, 7).Value = "Да"
$ws.Cells.Item(12, 8).Value = "Доволен"
$ws.Cells.Item(12, 9).Value = "Да"

# --- append new row 13 -------------------------------------------------------
$ws.Cells.Item(13, 1).Value = 22
$ws.Cells.Item(13, 2).Value = "м"
$ws.Cells.Item(13, 3).Value = "Жаворонок"
$ws.Cells.Item(13, 4).Value = 6
$ws.Cells.Item(13, 5).Value = "Да"
$ws.Cells.Item(13, 6).Value = "Интроверт"
$ws.Cells.Item(13, 7).Value = "Да"
$ws.Cells.Item(13, 8).Value = "Доволен"
$ws.Cells.Item(13, 9).Value = "Нет"

# --- move the active selection the way the author left it ------------------
$ws.Range("E12").Select()
